# New weekly data for "Terminal Hortofrutícola Agro Chillán - Pera" workbook.
# Two new price records (fecha = 45093) are added at the top of the data
# block (rows 294-295), pushing all existing data rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 294:295 - Excel shifts rows 294:377 down to 296:379
# and the freshly inserted rows inherit the formatting of the row above
# (in particular column D keeps its date number format).
$ws.Rows("294:295").Insert()

# Row 294: Packham's Triumph - Especial
$ws.Range("A294").Value = 7
$ws.Range("B294").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C294").Value = "Ñuble"
$ws.Range("D294").Value = 45093
$ws.Range("E294").Value = 16
$ws.Range("F294").Value = "Fruta"
$ws.Range("G294").Value = 100104
$ws.Range("H294").Value = "Frutos de pepita"
$ws.Range("I294").Value = 100104005
$ws.Range("J294").Value = "Pera"
$ws.Range("K294").Value = "Packham's Triumph"
$ws.Range("L294").Value = "Especial"
$ws.Range("M294").Value = 120
$ws.Range("N294").Value = 12000
$ws.Range("O294").Value = 12000
$ws.Range("P294").Value = 12000
$ws.Range("Q294").Value = "$/bandeja 18 kilos granel"
$ws.Range("R294").Value = "Región de O'Higgins"
$ws.Range("S294").Value = 667
$ws.Range("T294").Value = 18

# Row 295: Packham's Triumph - Primera
$ws.Range("A295").Value = 7
$ws.Range("B295").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C295").Value = "Ñuble"
$ws.Range("D295").Value = 45093
$ws.Range("E295").Value = 16
$ws.Range("F295").Value = "Fruta"
$ws.Range("G295").Value = 100104
$ws.Range("H295").Value = "Frutos de pepita"
$ws.Range("I295").Value = 100104005
$ws.Range("J295").Value = "Pera"
$ws.Range("K295").Value = "Packham's Triumph"
$ws.Range("L295").Value = "Primera"
$ws.Range("M295").Value = 100
$ws.Range("N295").Value = 10000
$ws.Range("O295").Value = 10000
$ws.Range("P295").Value = 10000
$ws.Range("Q295").Value = "$/bandeja 18 kilos granel"
$ws.Range("R295").Value = "Región de O'Higgins"
$ws.Range("S295").Value = 556
$ws.Range("T295").Value = 18
